$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 425.55554
$ws.Range("I2").Value = 425.55554
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 425.55554
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -312.55554
$ws.Range("N2").ClearContents()
# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 113
$ws.Range("N4").Value = -231
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 909116.2
$ws.Range("I6").Value = 1111132.6
$ws.Range("K6").Value = 3333397.8
$ws.Range("M6").Value = -3333285.8
# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 22.5
$ws.Range("I8").Value = 22.5
$ws.Range("K8").Value = 67.5
$ws.Range("M8").Value = 71.5
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 773302.1
$ws.Range("J17").Value = 773302.1
$ws.Range("L17").Value = 2319906.3
$ws.Range("N17").Value = -2320242.3
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2467.7778
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2467.7778
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2467.7778
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2817.7778
# Row 52 (Leve Item ID 4567)
$ws.Range("H52").Value = 1786
$ws.Range("I52").Value = 350
$ws.Range("J52").Value = 2743.3333
$ws.Range("K52").Value = 1050
$ws.Range("L52").Value = 8229.999899999999
$ws.Range("M52").Value = -890
$ws.Range("N52").Value = -8549.999899999999
# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 25129832
$ws.Range("I62").Value = 25129832
$ws.Range("K62").Value = 25129832
$ws.Range("M62").Value = -25129208
# Row 63 (Leve Item ID 10652)
$ws.Range("H63").Value = 89704.25
$ws.Range("J63").Value = 89704.25
$ws.Range("L63").Value = 89704.25
$ws.Range("N63").Value = -90952.25
# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 25129832
$ws.Range("I65").Value = 25129832
$ws.Range("K65").Value = 125649160
$ws.Range("M65").Value = -125646040
# Row 66 (Leve Item ID 10652)
$ws.Range("H66").Value = 89704.25
$ws.Range("J66").Value = 89704.25
$ws.Range("L66").Value = 269112.75
$ws.Range("N66").Value = -275352.75
# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 6467
$ws.Range("J70").Value = 8749.5
$ws.Range("L70").Value = 26248.5
$ws.Range("N70").Value = -26788.5
# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 6467
$ws.Range("J73").Value = 8749.5
$ws.Range("L73").Value = 26248.5
$ws.Range("N73").Value = -28120.5
# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 5354.143
$ws.Range("I86").Value = 4597
$ws.Range("K86").Value = 4597
$ws.Range("M86").Value = -3474
# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 5354.143
$ws.Range("I89").Value = 4597
$ws.Range("K89").Value = 22985
$ws.Range("M89").Value = -17369
# Row 103 (Leve Item ID 19909)
$ws.Range("H103").Value = 713
$ws.Range("J103").Value = 713
$ws.Range("L103").Value = 2139
$ws.Range("N103").Value = -3311
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 419.95834
$ws.Range("I107").Value = 367.05264
$ws.Range("J107").Value = 621
$ws.Range("K107").Value = 367.05264
$ws.Range("L107").Value = 621
$ws.Range("M107").Value = 1552.94736
$ws.Range("N107").Value = -4461
# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 73270.86
$ws.Range("J112").Value = 85227.75
$ws.Range("L112").Value = 255683.25
$ws.Range("N112").Value = -257899.25
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 12581.167
$ws.Range("I116").Value = 15995.833
$ws.Range("K116").Value = 15995.833
$ws.Range("M116").Value = -12553.833
# Row 130 (Leve Item ID 34691)
$ws.Range("H130").Value = 130742
$ws.Range("J130").Value = 130742
$ws.Range("L130").Value = 130742
$ws.Range("N130").Value = -140782
# Row 131 (Leve Item ID 36108)
$ws.Range("H131").Value = 11817.3125
$ws.Range("I131").Value = 3264.8
$ws.Range("J131").Value = 140105
$ws.Range("K131").Value = 9794.400000000001
$ws.Range("L131").Value = 420315
$ws.Range("M131").Value = -4754.400000000001
$ws.Range("N131").Value = -430395
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2680.3157
$ws.Range("I132").Value = 2094.4614
$ws.Range("J132").Value = 3949.6667
$ws.Range("K132").Value = 6283.3842
$ws.Range("L132").Value = 11849.0001
$ws.Range("M132").Value = -3753.3842
$ws.Range("N132").Value = -16909.0001
# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 1108
$ws.Range("I135").Value = 774.5714
$ws.Range("J135").Value = 2275
$ws.Range("K135").Value = 6971.1426
$ws.Range("L135").Value = 20475
$ws.Range("M135").Value = -4436.1426
$ws.Range("N135").Value = -25545
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1050.9286
$ws.Range("I137").Value = 985.61536
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 2956.84608
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -406.8460800000003
$ws.Range("N137").Value = -10800
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3227.7068
$ws.Range("I138").Value = 1463.8667
$ws.Range("J138").Value = 4403.6
$ws.Range("K138").Value = 4391.6001
$ws.Range("L138").Value = 13210.8
$ws.Range("M138").Value = 748.3999000000003
$ws.Range("N138").Value = -23490.8
# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 2300
$ws.Range("I141").Value = 2300
$ws.Range("K141").Value = 6900
$ws.Range("M141").Value = -1720

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 91004050
$ws.Range("I2").Value = 125129064
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 125129064
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -125128951
$ws.Range("N2").Value = -4226
# Row 3 (Leve Item ID 2494)
$ws.Range("H3").Value = 40033.332
$ws.Range("J3").Value = 60000
$ws.Range("L3").Value = 60000
$ws.Range("N3").Value = -60230
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5042.1143
$ws.Range("I32").Value = 4745
$ws.Range("K32").Value = 4745
$ws.Range("M32").Value = -4458
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 102825.37
$ws.Range("I45").Value = 112411
$ws.Range("J45").Value = 6969
$ws.Range("K45").Value = 112411
$ws.Range("L45").Value = 6969
$ws.Range("M45").Value = -112034
$ws.Range("N45").Value = -7723
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2854.9019
$ws.Range("I61").Value = 2477.6
$ws.Range("J61").Value = 4226.909
$ws.Range("K61").Value = 2477.6
$ws.Range("L61").Value = 4226.909
$ws.Range("M61").Value = -2265.6
$ws.Range("N61").Value = -4650.909
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 10567.5
$ws.Range("I74").Value = 1583.7333
$ws.Range("K74").Value = 1583.7333
$ws.Range("M74").Value = -709.7333000000001
# Row 75 (Leve Item ID 10714)
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26748
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 10567.5
$ws.Range("I77").Value = 1583.7333
$ws.Range("K77").Value = 7918.6665
$ws.Range("M77").Value = -3550.6665
# Row 78 (Leve Item ID 10714)
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -83736
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 638.2727
$ws.Range("I97").Value = 602.3333
$ws.Range("K97").Value = 602.3333
$ws.Range("M97").Value = -106.3333
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 5985.3213
$ws.Range("I110").Value = 6656.95
$ws.Range("J110").Value = 4306.25
$ws.Range("K110").Value = 6656.95
$ws.Range("L110").Value = 4306.25
$ws.Range("M110").Value = -4611.95
$ws.Range("N110").Value = -8396.25
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 91004050
$ws.Range("I116").Value = 125129064
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 125129064
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -125126770
$ws.Range("N116").Value = -8588
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2930.5757
$ws.Range("I122").Value = 2855.6897
$ws.Range("J122").Value = 3473.5
$ws.Range("K122").Value = 8567.069100000001
$ws.Range("L122").Value = 10420.5
$ws.Range("M122").Value = -6117.069100000001
$ws.Range("N122").Value = -15320.5
# Row 131 (Leve Item ID 34706)
$ws.Range("H131").Value = 94995.8
$ws.Range("J131").Value = 94995.8
$ws.Range("L131").Value = 94995.8
$ws.Range("N131").Value = -105075.8
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3443.875
$ws.Range("I132").Value = 3011.9583
$ws.Range("K132").Value = 9035.874899999999
$ws.Range("M132").Value = -6505.874899999999
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2854.9019
$ws.Range("I136").Value = 2477.6
$ws.Range("J136").Value = 4226.909
$ws.Range("K136").Value = 7432.799999999999
$ws.Range("L136").Value = 12680.727
$ws.Range("M136").Value = -4882.799999999999
$ws.Range("N136").Value = -17780.727

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 91004050
$ws.Range("I3").Value = 125129064
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 125129064
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -125128950
$ws.Range("N3").Value = -4228
# Row 11 (Leve Item ID 2481)
$ws.Range("H11").Value = 244.85715
$ws.Range("I11").Value = 102.25
$ws.Range("J11").Value = 435
$ws.Range("K11").Value = 102.25
$ws.Range("L11").Value = 435
$ws.Range("M11").Value = 37.75
$ws.Range("N11").Value = -715
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 406525.75
$ws.Range("I22").Value = 672.3125
$ws.Range("K22").Value = 672.3125
$ws.Range("M22").Value = -499.3125
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1999.9231
$ws.Range("I86").Value = 1999.9231
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1999.9231
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -876.9231
$ws.Range("N86").ClearContents()
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1999.9231
$ws.Range("I89").Value = 1999.9231
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9999.6155
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4383.6155
$ws.Range("N89").ClearContents()
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 2900
$ws.Range("I94").Value = 2900
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2900
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2449
$ws.Range("N94").ClearContents()
# Row 100 (Leve Item ID 18347)
$ws.Range("H100").Value = 33250
$ws.Range("J100").Value = 33250
$ws.Range("L100").Value = 33250
$ws.Range("N100").Value = -35414
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 2762.6155
$ws.Range("I107").Value = 2139.7144
$ws.Range("J107").Value = 3489.3333
$ws.Range("K107").Value = 2139.7144
$ws.Range("L107").Value = 3489.3333
$ws.Range("M107").Value = -219.7143999999998
$ws.Range("N107").Value = -7329.3333
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1595.3064
$ws.Range("I134").Value = 1581.8167
$ws.Range("K134").Value = 4745.4501
$ws.Range("M134").Value = -2210.4501

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 13 (Leve Item ID 1996)
$ws.Range("H13").Value = 3812.5
$ws.Range("I13").Value = 4000
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = -3861
$ws.Range("N13").Value = -2778
# Row 21 (Leve Item ID 2000)
$ws.Range("H21").Value = 1000
$ws.Range("J21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("N21").Value = -1470
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 629
$ws.Range("I22").Value = 649
$ws.Range("J22").Value = 609
$ws.Range("K22").Value = 649
$ws.Range("L22").Value = 609
$ws.Range("M22").Value = -299
$ws.Range("N22").Value = -1309
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 27076.512
$ws.Range("I31").Value = 43657.082
$ws.Range("J31").Value = 3668.647
$ws.Range("K31").Value = 43657.082
$ws.Range("L31").Value = 3668.647
$ws.Range("M31").Value = -43362.082
$ws.Range("N31").Value = -4258.647
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 27076.512
$ws.Range("I34").Value = 43657.082
$ws.Range("J34").Value = 3668.647
$ws.Range("K34").Value = 43657.082
$ws.Range("L34").Value = 3668.647
$ws.Range("M34").Value = -43455.082
$ws.Range("N34").Value = -4072.647
# Row 57 (Leve Item ID 3908)
$ws.Range("H57").Value = 29999
$ws.Range("J57").Value = 29999
$ws.Range("L57").Value = 29999
$ws.Range("N57").Value = -31119
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 2284.1018
$ws.Range("I58").Value = 2288.025
$ws.Range("J58").Value = 2275.842
$ws.Range("K58").Value = 2288.025
$ws.Range("L58").Value = 2275.842
$ws.Range("M58").Value = -2085.025
$ws.Range("N58").Value = -2681.842
# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 26000
$ws.Range("J59").Value = 26000
$ws.Range("L59").Value = 26000
$ws.Range("N59").Value = -28290
# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748
# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -98736
# Row 87 (Leve Item ID 11929)
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90 (Leve Item ID 11929)
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 6504.8335
$ws.Range("I99").Value = 2577.7144
$ws.Range("K99").Value = 2577.7144
$ws.Range("M99").Value = -1079.7144
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1335
$ws.Range("I107").Value = 1034.4
$ws.Range("J107").Value = 1460.25
$ws.Range("K107").Value = 1034.4
$ws.Range("L107").Value = 1460.25
$ws.Range("M107").Value = 885.5999999999999
$ws.Range("N107").Value = -5300.25
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1654
$ws.Range("I122").Value = 1555.4
$ws.Range("K122").Value = 4666.200000000001
$ws.Range("M122").Value = -2216.200000000001
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 6504.8335
$ws.Range("I126").Value = 2577.7144
$ws.Range("K126").Value = 7733.1432
$ws.Range("M126").Value = -5263.1432
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 4625.491
$ws.Range("I132").Value = 3154.449
$ws.Range("J132").Value = 13635.625
$ws.Range("K132").Value = 9463.347
$ws.Range("L132").Value = 40906.875
$ws.Range("M132").Value = -6933.347
$ws.Range("N132").Value = -45966.875
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 13856.368
$ws.Range("I134").Value = 7524.4546
$ws.Range("J134").Value = 55647
$ws.Range("K134").Value = 22573.3638
$ws.Range("L134").Value = 166941
$ws.Range("M134").Value = -20038.3638
$ws.Range("N134").Value = -172011
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 2284.1018
$ws.Range("I136").Value = 2288.025
$ws.Range("J136").Value = 2275.842
$ws.Range("K136").Value = 6864.075000000001
$ws.Range("L136").Value = 6827.526
$ws.Range("M136").Value = -4314.075000000001
$ws.Range("N136").Value = -11927.526

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 352.5
$ws.Range("J12").Value = 476.9
$ws.Range("L12").Value = 1430.7
$ws.Range("N12").Value = -1776.7
# Row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 884.2857
$ws.Range("I23").Value = 2750
$ws.Range("J23").Value = 573.3333
$ws.Range("K23").Value = 8250
$ws.Range("L23").Value = 1719.9999
$ws.Range("M23").Value = -8015
$ws.Range("N23").Value = -2189.9999
# Row 50 (Leve Item ID 4725)
$ws.Range("H50").Value = 1590.2142
$ws.Range("J50").Value = 3116.6667
$ws.Range("L50").Value = 9350.000100000001
$ws.Range("N50").Value = -10312.0001
# Row 53 (Leve Item ID 4725)
$ws.Range("H53").Value = 1590.2142
$ws.Range("J53").Value = 3116.6667
$ws.Range("L53").Value = 9350.000100000001
$ws.Range("N53").Value = -10312.0001
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 915.6429000000001
$ws.Range("I107").Value = 377
$ws.Range("K107").Value = 1131
$ws.Range("M107").Value = 789
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1096.4572
$ws.Range("I132").Value = 1072.303
$ws.Range("J132").Value = 1495
$ws.Range("K132").Value = 9650.727000000001
$ws.Range("L132").Value = 13455
$ws.Range("M132").Value = -7120.727000000001
$ws.Range("N132").Value = -18515
# Row 141 (Leve Item ID 44076)
$ws.Range("H141").Value = 206712.6
$ws.Range("I141").Value = 8382.5
$ws.Range("K141").Value = 25147.5
$ws.Range("M141").Value = -19967.5

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 12 (Leve Item ID 4093)
$ws.Range("H12").Value = 2000
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3280
# Row 18 (Leve Item ID 4309)
$ws.Range("H18").Value = 37072036
$ws.Range("I18").Value = 111111110
$ws.Range("K18").Value = 111111110
$ws.Range("M18").Value = -111110817
# Row 20 (Leve Item ID 4095)
$ws.Range("H20").Value = 7919.6665
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 24 (Leve Item ID 4431)
$ws.Range("H24").Value = 80000
$ws.Range("J24").Value = 80000
$ws.Range("L24").Value = 80000
$ws.Range("N24").Value = -80346
# Row 45 (Leve Item ID 27225)
$ws.Range("H45").Value = 27450
$ws.Range("J45").Value = 27450
$ws.Range("L45").Value = 27450
$ws.Range("N45").Value = -28568
# Row 51 (Leve Item ID 27222)
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26018
# Row 55 (Leve Item ID 4237)
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 13210.2
$ws.Range("I57").Value = 6055
$ws.Range("K57").Value = 6055
$ws.Range("M57").Value = -5235
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 28723.2
$ws.Range("I70").Value = 23999
$ws.Range("K70").Value = 23999
$ws.Range("M70").Value = -23729
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 28723.2
$ws.Range("I73").Value = 23999
$ws.Range("K73").Value = 23999
$ws.Range("M73").Value = -23063
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2960.75
$ws.Range("I80").Value = 2960.75
$ws.Range("K80").Value = 2960.75
$ws.Range("M80").Value = -1962.75
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2960.75
$ws.Range("I83").Value = 2960.75
$ws.Range("K83").Value = 14803.75
$ws.Range("M83").Value = -9811.75
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 76925590
$ws.Range("I102").Value = 2184.5557
$ws.Range("K102").Value = 2184.5557
$ws.Range("M102").Value = -562.5556999999999
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 637.13336
$ws.Range("I107").Value = 614.5454999999999
$ws.Range("J107").Value = 699.25
$ws.Range("K107").Value = 614.5454999999999
$ws.Range("L107").Value = 699.25
$ws.Range("M107").Value = 1305.4545
$ws.Range("N107").Value = -4539.25
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1583.3
$ws.Range("I122").Value = 1205.8823
$ws.Range("K122").Value = 3617.6469
$ws.Range("M122").Value = -1167.6469
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3838.8635
$ws.Range("I132").Value = 3549.5789
$ws.Range("K132").Value = 10648.7367
$ws.Range("M132").Value = -8118.736699999999

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 20 (Leve Item ID 4308)
$ws.Range("H20").Value = 80000
$ws.Range("J20").Value = 80000
$ws.Range("L20").Value = 80000
$ws.Range("N20").Value = -80452
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2195.6956
$ws.Range("I22").Value = 1785.8572
$ws.Range("J22").Value = 2375
$ws.Range("K22").Value = 1785.8572
$ws.Range("L22").Value = 2375
$ws.Range("M22").Value = -1490.8572
$ws.Range("N22").Value = -2965
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2195.6956
$ws.Range("I27").Value = 1785.8572
$ws.Range("J27").Value = 2375
$ws.Range("K27").Value = 1785.8572
$ws.Range("L27").Value = 2375
$ws.Range("M27").Value = -1678.8572
$ws.Range("N27").Value = -2589
# Row 29 (Leve Item ID 3576)
$ws.Range("H29").Value = 3807.8333
$ws.Range("J29").Value = 3769.4
$ws.Range("L29").Value = 3769.4
$ws.Range("N29").Value = -4359.4
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 4833.6665
$ws.Range("I40").Value = 4647.4116
$ws.Range("K40").Value = 4647.4116
$ws.Range("M40").Value = -4511.4116
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 806.61536
$ws.Range("I46").Value = 707.1667
$ws.Range("K46").Value = 707.1667
$ws.Range("M46").Value = -519.1667
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 816.6875
$ws.Range("I55").Value = 642.7143
$ws.Range("J55").Value = 952
$ws.Range("K55").Value = 642.7143
$ws.Range("L55").Value = 952
$ws.Range("M55").Value = -469.7143
$ws.Range("N55").Value = -1298
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3072
$ws.Range("I68").Value = 2441.6667
$ws.Range("K68").Value = 2441.6667
$ws.Range("M68").Value = -1692.6667
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3072
$ws.Range("I71").Value = 2441.6667
$ws.Range("K71").Value = 12208.3335
$ws.Range("M71").Value = -8464.333500000001
# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 1348.697
$ws.Range("I82").Value = 1256.64
$ws.Range("K82").Value = 1256.64
$ws.Range("M82").Value = -895.6400000000001
# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 1348.697
$ws.Range("I85").Value = 1256.64
$ws.Range("K85").Value = 1256.64
$ws.Range("M85").Value = -8.6400000000001
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 345
$ws.Range("I93").Value = 345
$ws.Range("K93").Value = 345
$ws.Range("M93").Value = 903
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 1661.5
$ws.Range("I100").Value = 1593.8
$ws.Range("K100").Value = 1593.8
$ws.Range("M100").Value = -1052.8
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5428.9546
$ws.Range("I132").Value = 4829.8667
$ws.Range("K132").Value = 14489.6001
$ws.Range("M132").Value = -11959.6001
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 5004.3184
$ws.Range("I136").Value = 4397.467
$ws.Range("K136").Value = 13192.401
$ws.Range("M136").Value = -10642.401

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 49737.383
$ws.Range("J4").Value = 2515.0715
$ws.Range("L4").Value = 2515.0715
$ws.Range("N4").Value = -2741.0715
# Row 20 (Leve Item ID 3023)
$ws.Range("H20").Value = 2005877.6
$ws.Range("I20").Value = 10000000
$ws.Range("J20").Value = 7347
$ws.Range("K20").Value = 10000000
$ws.Range("L20").Value = 7347
$ws.Range("M20").Value = -9999760
$ws.Range("N20").Value = -7827
# Row 26 (Leve Item ID 3800)
$ws.Range("H26").Value = 15202.4
$ws.Range("I26").Value = 9670.666999999999
$ws.Range("K26").Value = 9670.666999999999
$ws.Range("M26").Value = -9377.666999999999
# Row 28 (Leve Item ID 3053)
$ws.Range("H28").Value = 80000
$ws.Range("J28").Value = 80000
$ws.Range("L28").Value = 80000
$ws.Range("N28").Value = -80696
# Row 31 (Leve Item ID 3052)
$ws.Range("H31").Value = 80000
$ws.Range("J31").Value = 80000
$ws.Range("L31").Value = 80000
$ws.Range("N31").Value = -80696
# Row 43 (Leve Item ID 3831)
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 10290.765
$ws.Range("I81").Value = 34147.668
$ws.Range("J81").Value = 5178.5713
$ws.Range("K81").Value = 68295.336
$ws.Range("L81").Value = 10357.1426
$ws.Range("M81").Value = -67234.336
$ws.Range("N81").Value = -12479.1426
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 10290.765
$ws.Range("I84").Value = 34147.668
$ws.Range("J84").Value = 5178.5713
$ws.Range("K84").Value = 341476.68
$ws.Range("L84").Value = 51785.713
$ws.Range("M84").Value = -336172.68
$ws.Range("N84").Value = -62393.713
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 2407
$ws.Range("I96").Value = 2444.6
$ws.Range("K96").Value = 2444.6
$ws.Range("M96").Value = -1071.6
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 1638.2
$ws.Range("I113").Value = 1554.2
$ws.Range("K113").Value = 4662.6
$ws.Range("M113").Value = -2492.6
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2438.375
$ws.Range("I122").Value = 2294.0286
$ws.Range("K122").Value = 6882.085800000001
$ws.Range("M122").Value = -4432.085800000001
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 3880
$ws.Range("I126").Value = 3100
$ws.Range("K126").Value = 9300
$ws.Range("M126").Value = -6830
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1504.1351
$ws.Range("I132").Value = 1568.4333
$ws.Range("J132").Value = 1228.5714
$ws.Range("K132").Value = 4705.2999
$ws.Range("L132").Value = 3685.7142
$ws.Range("M132").Value = -2175.2999
$ws.Range("N132").Value = -8745.7142
# Row 133 (Leve Item ID 41869)
$ws.Range("H133").Value = 26440
$ws.Range("J133").Value = 25550.334
$ws.Range("L133").Value = 25550.334
$ws.Range("N133").Value = -35670.334
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2114.1428
$ws.Range("I136").Value = 1653.6538
$ws.Range("J136").Value = 3444.4443
$ws.Range("K136").Value = 4960.9614
$ws.Range("L136").Value = 10333.3329
$ws.Range("M136").Value = -2410.9614
$ws.Range("N136").Value = -15433.3329
